$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.369.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = "'1.842.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("D4").Value = "'0.9989"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'238.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").Value = "'0.6303"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.07527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("D9").Value = "'0.2926"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("D10").Value = "'24.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").Value = "'1.866.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.98%  '
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("D15").Value = "'0.00001030"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.06%  '
$ws.Range("D16").Value = "'82.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = "'2.119.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.39%  '
$ws.Range("D18").Value = "'6.144"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").Value = "'29.410.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = "'227.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.73%  '
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").Value = "'7.435"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.78%  '
$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").Value = "'156.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("D26").Value = "'0.1389"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'8.352"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.03%  '
$ws.Range("D28").Value = "'17.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.40%  '
$ws.Range("D29").Value = "'1.456"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").Value = "'1.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.21%  '
$ws.Range("D31").Value = "'0.05621"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("D33").Value = "'4.018"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("D34").Value = "'1.831"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.14%  '
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").Value = "'0.7112"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("D37").Value = "'2.588"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").Value = "'1.241.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("D39").Value = "'0.01804"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("D40").Value = "'2.760"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").Value = "'6.324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.80%  '
$ws.Range("D42").Value = "'0.9005"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("D43").Value = "'0.9993"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'101.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("E45").Value = '  -2.09%  '
$ws.Range("D46").Value = "'0.00000000120"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.62%  '
$ws.Range("D47").Value = "'7.056"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.65%  '
$ws.Range("E48").Value = '  -0.49%  '
$ws.Range("E49").Value = '  -1.69%  '
$ws.Range("D50").Value = "'8.869"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.15%  '
$ws.Range("E51").Value = '  -0.75%  '
